# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" sheet, formatted
#    like the other quarterly sheets, holding the new quarter's fund holdings
#    detail.
# 2. Insert a new top data row into "总计" summarizing the new quarter,
#    pushing the older rows down one row (and renumbering the index column).

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")
$templateSheet = $wb.Worksheets.Item("2021-Q4")

$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# The worksheet collection re-indexes on insert, so re-resolve "总计" by name
# now that the new sheet has been placed in front of it.
$totalSheet = $wb.Worksheets.Item("总计")

# Copy header + index-column formatting from an existing quarterly sheet so the
# new sheet matches the established look (bold centered header, bordered A col).
$templateSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$templateSheet.Range("A2:A9").Copy()
$newSheet.Range("A2:A9").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Fund code / name / size / position columns are stored as plain text in this
# workbook (leading zeros in fund codes, fixed decimal strings like "6.20"
# must survive verbatim) -- format as Text before writing so COM doesn't
# coerce them to numbers.
$newSheet.Range("B2:G9").NumberFormat = "@"

$fundRows = @(
    @(0, "010336", "中欧悦享生活混合A", "44.28", "90.44", "6.08", "2.6922", 7),
    @(1, "002621", "中欧消费主题股票A", "19.29", "88.29", "7.05", "1.3599", 5),
    @(2, "002697", "中欧消费主题股票C", "6.20", "88.29", "7.05", "0.4371", 5),
    @(3, "005620", "中欧品质消费股票A", "3.74", "90.47", "7.16", "0.2678", 6),
    @(4, "010852", "中欧内需成长混合型证券投资基金A", "5.23", "91.46", "5.11", "0.2673", 7),
    @(5, "005621", "中欧品质消费股票C", "1.11", "90.47", "7.16", "0.0795", 6),
    @(6, "010337", "中欧悦享生活混合C", "1.08", "90.44", "6.08", "0.0657", 7),
    @(7, "010853", "中欧内需成长混合型证券投资基金C", "0.67", "91.46", "5.11", "0.0342", 7)
)

$r = 2
foreach ($row in $fundRows) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $newSheet.Cells.Item($r, 5).Value = $row[4]
    $newSheet.Cells.Item($r, 6).Value = $row[5]
    $newSheet.Cells.Item($r, 7).Value = $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Shift the "总计" rows down by one to make room for the new 2022-Q1 row, then
# write the new row at the top and renumber the index column.
$totalSheet.Range("A2:D5").Copy()
$totalSheet.Range("A3:D6").PasteSpecial(-4104)

# The paste above extends the sheet past its previous last row (row 5 -> 6);
# the new row's A6 cell doesn't pick up the source formatting from the
# all-in-one paste, so re-apply just the style for that cell explicitly.
$totalSheet.Range("A5").Copy()
$totalSheet.Range("A6").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 8
$totalSheet.Range("D2").Value = 5.2

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
